# Order collection api in progress
# Update existing order rows (2-4) with new data and append two new order
# rows (5-6) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force Excel to store the value as literal text (matching shared-string
    # cells for numeric-looking data like phone numbers / zip codes /
    # amounts) rather than letting a bare Value assignment auto-convert it
    # to a number. Using a quoted formula then flattening it via
    # copy/paste-values keeps the cell's style untouched (no NumberFormat
    # change survives on the cell).
    $escaped = $value.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# --- Row 2 : KEYCLUE order for 樊聪 (fancong6985) -----------------------
$ws.Range("A2").Value = "2018-04-18 22:20:38"
$ws.Range("B2").Value = "KEYCLUE海外旗舰店"
$ws.Range("C2").Value = "fancong6985"
$ws.Range("D2").Value = "2018-04-18 22:20:42"
$ws.Range("E2").Value = "樊聪"
Set-TextValue $ws.Range("F2") "13811765578"
Set-TextValue $ws.Range("G2") "100011"
$ws.Range("H2").Value = "北京"
$ws.Range("I2").Value = "北京市"
$ws.Range("J2").Value = "朝阳区"
$ws.Range("K2").Value = "奥运村街道朝阳区域清街2号院融域嘉园4号楼6单元801"
$ws.Range("L2").Value = "EEOG1WSR02W_YLF"
Set-TextValue $ws.Range("M2") "1241.04"
$ws.Range("N2").ClearContents()
$ws.Range("O2").Value = "eyeye春季新款韩版潮流时尚V领条纹蝴蝶袖上衣T恤"

# --- Row 3 : KEYCLUE order for 郭舒君 (郭舒君1234) ------------------------
$ws.Range("A3").Value = "2018-04-18 20:44:27"
$ws.Range("B3").Value = "KEYCLUE海外旗舰店"
$ws.Range("C3").Value = "郭舒君1234"
$ws.Range("D3").Value = "2018-04-18 20:44:32"
$ws.Range("E3").Value = "郭舒君"
Set-TextValue $ws.Range("F3") "13829910807"
Set-TextValue $ws.Range("G3") "516001"
$ws.Range("H3").Value = "广东省"
$ws.Range("I3").Value = "惠州市"
$ws.Range("J3").Value = "惠城区"
$ws.Range("K3").Value = "江北街道江北佳兆业二期A3座806房"
$ws.Range("L3").Value = "mnmm17ssc21_7p"
Set-TextValue $ws.Range("M3") "323.39"
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = "MIDNIGHT MOMENT.人造珍珠链条iPhone7/6苹果手机保护壳韩国"

# --- Row 4 : KEYCLUE order for 汪凯旋 (凯凯99977573209) -------------------
$ws.Range("A4").Value = "2018-04-17 22:23:42"
$ws.Range("B4").Value = "KEYCLUE海外旗舰店"
$ws.Range("C4").Value = "凯凯99977573209"
$ws.Range("D4").Value = "2018-04-17 22:25:08"
$ws.Range("E4").Value = "汪凯旋"
Set-TextValue $ws.Range("F4") "13640787455"
Set-TextValue $ws.Range("G4") "510180"
$ws.Range("H4").Value = "广东省"
$ws.Range("I4").Value = "广州市"
$ws.Range("J4").Value = "越秀区"
$ws.Range("K4").Value = "珠光街道 沿江中路299号海俊酒店前台"
$ws.Range("L4").Value = "S72BL08BK"
Set-TextValue $ws.Range("M4") "399.00"
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = "Salad bowls官方正品韩版甜美公主风一字肩短袖韩国包邮"

# --- Row 5 : KEYCLUE order for 楼旦韵 (takigo92) [new row] ---------------
$ws.Range("A5").Value = "2018-04-17 20:38:01"
$ws.Range("B5").Value = "KEYCLUE海外旗舰店"
$ws.Range("C5").Value = "takigo92"
$ws.Range("D5").Value = "2018-04-17 20:38:05"
$ws.Range("E5").Value = "楼旦韵"
Set-TextValue $ws.Range("F5") "13616511558"
Set-TextValue $ws.Range("G5") "310013"
$ws.Range("H5").Value = "浙江省"
$ws.Range("I5").Value = "杭州市"
$ws.Range("J5").Value = "西湖区"
$ws.Range("K5").Value = "杭州市西湖区黄姑山路38-1号十足(菜鸟驿站:0571-88116223)"
$ws.Range("L5").Value = "S72TS50IV"
Set-TextValue $ws.Range("M5") "299.00"
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = "Salad bowls官方正品韩版可爱清新假V领百搭短袖韩国包邮"

# --- Row 6 : KEYCLUE order for 泛泛 (missvan82) [new row] ----------------
$ws.Range("A6").Value = "2018-04-12 22:43:26"
$ws.Range("B6").Value = "KEYCLUE海外旗舰店"
$ws.Range("C6").Value = "missvan82"
$ws.Range("D6").Value = "2018-04-12 22:43:36"
$ws.Range("E6").Value = "泛泛"
Set-TextValue $ws.Range("F6") "18667162213"
Set-TextValue $ws.Range("G6") "310012"
$ws.Range("H6").Value = "浙江省"
$ws.Range("I6").Value = "杭州市"
$ws.Range("J6").Value = "西湖区"
$ws.Range("K6").Value = "三墩镇白马尊邸11-2-301"
$ws.Range("L6").Value = "S72TS12NV"
Set-TextValue $ws.Range("M6") "279.00"
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = "Salad bowls官方正品韩版宽松宽条纹圆领彩色短袖韩国包邮"
